$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the worked hours on row 10, Friday (column F)
$ws.Range("F10").Value = 6.25

# Update the selected cell to reflect where the user clicked last
$ws.Range("G16").Select()

